# Recalculation of quality indicators: renumber/rename rule & indicator
# codes on several sheets, drop an obsolete measure row, and add four new
# metric rows.
#
# Helper: write a value as genuine TEXT (so strings that look numeric, like
# "94.74%", are not silently reinterpreted by Excel as a percentage number)
# without leaving any stray number-format style behind on the cell.
function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) SCHEME_MEASURES: rename indicator codes MQMSxx -> MQME00x
# ---------------------------------------------------------------------------
$wsMeasures = $wb.Worksheets.Item("SCHEME_MEASURES")
Set-TextValue $wsMeasures.Range("A2") "MQME001"
Set-TextValue $wsMeasures.Range("A3") "MQME002"
Set-TextValue $wsMeasures.Range("A4") "MQME003"
Set-TextValue $wsMeasures.Range("A5") "MQME004"
Set-TextValue $wsMeasures.Range("A6") "MQME005"

# ---------------------------------------------------------------------------
# 2) METADATA_ISSUES: remap rule codes in column A (rows 2-78)
# ---------------------------------------------------------------------------
$wsIssues = $wb.Worksheets.Item("METADATA_ISSUES")

for ($r = 2; $r -le 3; $r++) {
    Set-TextValue $wsIssues.Cells.Item($r, 1) "MQME012"
}
for ($r = 4; $r -le 6; $r++) {
    Set-TextValue $wsIssues.Cells.Item($r, 1) "MQME014"
}
for ($r = 7; $r -le 72; $r++) {
    Set-TextValue $wsIssues.Cells.Item($r, 1) "MQME008"
}
for ($r = 73; $r -le 77; $r++) {
    Set-TextValue $wsIssues.Cells.Item($r, 1) "MQME009"
}
Set-TextValue $wsIssues.Cells.Item(78, 1) "MQME010"

# ---------------------------------------------------------------------------
# 3) METADATA_MEASURES: drop the "Total number of columns" row and renumber
#    the remaining two rows (MQMEA1 -> MQME006, MQMEA2 -> MQME007)
# ---------------------------------------------------------------------------
$wsMetaMeasures = $wb.Worksheets.Item("METADATA_MEASURES")
$wsMetaMeasures.Rows.Item(2).Delete()
Set-TextValue $wsMetaMeasures.Range("A2") "MQME006"
Set-TextValue $wsMetaMeasures.Range("A3") "MQME007"

# ---------------------------------------------------------------------------
# 4) METADATA_METRICS: renumber/relabel the existing rows and insert 4 new
#    metric rows (MQID008-MQID011) after the previously last row
# ---------------------------------------------------------------------------
$wsMetrics = $wb.Worksheets.Item("METADATA_METRICS")

Set-TextValue $wsMetrics.Range("A2") "MQID001"
Set-TextValue $wsMetrics.Range("B2") "Table names in singular"
Set-TextValue $wsMetrics.Range("C2") "94.74%"

Set-TextValue $wsMetrics.Range("A3") "MQID002"
Set-TextValue $wsMetrics.Range("B3") "Table with recommended name length"
Set-TextValue $wsMetrics.Range("C3") "100.00%"

Set-TextValue $wsMetrics.Range("A4") "MQID003"
Set-TextValue $wsMetrics.Range("B4") "Columns with correct prefixes"
Set-TextValue $wsMetrics.Range("C4") "99.22%"

Set-TextValue $wsMetrics.Range("A5") "MQID004"
Set-TextValue $wsMetrics.Range("B5") "Columns with recommended name size"
Set-TextValue $wsMetrics.Range("C5") "100.00%"

Set-TextValue $wsMetrics.Range("A6") "MQID005"
Set-TextValue $wsMetrics.Range("B6") "Columns with comments"
Set-TextValue $wsMetrics.Range("C6") "82.90%"

Set-TextValue $wsMetrics.Range("A7") "MQID006"
Set-TextValue $wsMetrics.Range("B7") "Table with standard PK prefixes"
Set-TextValue $wsMetrics.Range("C7") "84.38%"

Set-TextValue $wsMetrics.Range("A8") "MQID007"
Set-TextValue $wsMetrics.Range("B8") "Table with standard FK prefixes"
Set-TextValue $wsMetrics.Range("C8") "94.12%"

$wsMetrics.Rows.Item(9).Insert()
Set-TextValue $wsMetrics.Range("A9") "MQID008"
Set-TextValue $wsMetrics.Range("B9") "Table with standard UK prefixes"
Set-TextValue $wsMetrics.Range("C9") "0.00%"

$wsMetrics.Rows.Item(10).Insert()
Set-TextValue $wsMetrics.Range("A10") "MQID009"
Set-TextValue $wsMetrics.Range("B10") "NUMBER columns with valid scale"
Set-TextValue $wsMetrics.Range("C10") "100.00%"

$wsMetrics.Rows.Item(11).Insert()
Set-TextValue $wsMetrics.Range("A11") "MQID010"
Set-TextValue $wsMetrics.Range("B11") "Columns with valid num_distinct"
Set-TextValue $wsMetrics.Range("C11") "100.00%"

$wsMetrics.Rows.Item(12).Insert()
Set-TextValue $wsMetrics.Range("A12") "MQID011"
Set-TextValue $wsMetrics.Range("B12") "Columns with valid num_nulls"
Set-TextValue $wsMetrics.Range("C12") "100.00%"
